# Refresh cached FFXIV Tonberry-server market data across all Job sheets.
# Columns H:N on each sheet are:
#   H currentAveragePrice, I currentAveragePriceNQ, J currentAveragePriceHQ,
#   K LevePriceNQ, L LevePriceHQ, M LeveProfitNQ, N LeveProfitHQ
# This mirrors the scheduled runner that re-pulls Universalis prices and
# rewrites the affected cells' cached values in place (no formulas involved).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 2000
$ws.Range("I20").Value = 2000
$ws.Range("K20").Value = 2000
$ws.Range("M20").Value = -1770
$ws.Range("H35").Value = 2000
$ws.Range("I35").Value = 2000
$ws.Range("K35").Value = 2000
$ws.Range("M35").Value = -1621
$ws.Range("H125").Value = 1264.1666
$ws.Range("I125").Value = 1303.625
$ws.Range("J125").Value = 1185.25
$ws.Range("K125").Value = 11732.625
$ws.Range("L125").Value = 10667.25
$ws.Range("M125").Value = -9272.625
$ws.Range("N125").Value = -15587.25
$ws.Range("H132").Value = 887.44446
$ws.Range("I132").Value = 817.0513
$ws.Range("K132").Value = 2451.1539
$ws.Range("M132").Value = 78.84610000000021
$ws.Range("H133").Value = 89000
$ws.Range("J133").Value = 89000
$ws.Range("L133").Value = 89000
$ws.Range("N133").Value = -99120
$ws.Range("I137").Value = 1643.8334
$ws.Range("J137").Value = 2364.5715
$ws.Range("K137").Value = 4931.5002
$ws.Range("L137").Value = 7093.7145
$ws.Range("M137").Value = -2381.5002
$ws.Range("N137").Value = -12193.7145
$ws.Range("H138").Value = 1758.6558
$ws.Range("J138").Value = 2174.394
$ws.Range("L138").Value = 6523.181999999999
$ws.Range("N138").Value = -16803.182

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4469.174
$ws.Range("I61").Value = 3568.5
$ws.Range("K61").Value = 3568.5
$ws.Range("M61").Value = -3356.5
$ws.Range("H74").Value = 2905.25
$ws.Range("J74").Value = 3974.25
$ws.Range("L74").Value = 3974.25
$ws.Range("N74").Value = -5722.25
$ws.Range("H77").Value = 2905.25
$ws.Range("J77").Value = 3974.25
$ws.Range("L77").Value = 19871.25
$ws.Range("N77").Value = -28607.25
$ws.Range("H132").Value = 1558.2444
$ws.Range("I132").Value = 1059.7812
$ws.Range("K132").Value = 3179.3436
$ws.Range("M132").Value = -649.3435999999997
$ws.Range("H136").Value = 4469.174
$ws.Range("I136").Value = 3568.5
$ws.Range("K136").Value = 10705.5
$ws.Range("M136").Value = -8155.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5578.3335
$ws.Range("I134").Value = 6528.5386
$ws.Range("K134").Value = 19585.6158
$ws.Range("M134").Value = -17050.6158

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2161.862
$ws.Range("I31").Value = 1744.25
$ws.Range("K31").Value = 1744.25
$ws.Range("M31").Value = -1449.25
$ws.Range("H34").Value = 2161.862
$ws.Range("I34").Value = 1744.25
$ws.Range("K34").Value = 1744.25
$ws.Range("M34").Value = -1542.25
$ws.Range("H132").Value = 2431.8386
$ws.Range("I132").Value = 1599
$ws.Range("J132").Value = 3946.0908
$ws.Range("K132").Value = 4797
$ws.Range("L132").Value = 11838.2724
$ws.Range("M132").Value = -2267
$ws.Range("N132").Value = -16898.2724
$ws.Range("H134").Value = 1111.7046
$ws.Range("I134").Value = 1108.3077
$ws.Range("K134").Value = 3324.9231
$ws.Range("M134").Value = -789.9231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 841.92
$ws.Range("J131").Value = 853.3158
$ws.Range("L131").Value = 2559.9474
$ws.Range("N131").Value = -12639.9474
$ws.Range("H132").Value = 1268
$ws.Range("I132").Value = 890
$ws.Range("J132").Value = 1362.5
$ws.Range("K132").Value = 8010
$ws.Range("L132").Value = 12262.5
$ws.Range("M132").Value = -5480
$ws.Range("N132").Value = -17322.5
$ws.Range("H133").Value = 5000
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 29.333334
$ws.Range("J2").Value = 126.5
$ws.Range("K2").Value = 29.333334
$ws.Range("L2").Value = 126.5
$ws.Range("M2").Value = 83.66666599999999
$ws.Range("N2").Value = -352.5
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").ClearContents()
$ws.Range("H43").Value = 7006
$ws.Range("I43").Value = 7006
$ws.Range("K43").Value = 7006
$ws.Range("M43").Value = -6855
$ws.Range("H126").Value = 9429161
$ws.Range("I126").Value = 13892492
$ws.Range("J126").Value = 502500
$ws.Range("K126").Value = 41677476
$ws.Range("L126").Value = 1507500
$ws.Range("M126").Value = -41675006
$ws.Range("N126").Value = -1512440
$ws.Range("H132").Value = 2266647.5
$ws.Range("I132").Value = 3208684.2
$ws.Range("J132").Value = 5759.6
$ws.Range("K132").Value = 9626052.600000001
$ws.Range("L132").Value = 17278.8
$ws.Range("M132").Value = -9623522.600000001
$ws.Range("N132").Value = -22338.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3593.5833
$ws.Range("I7").Value = 2313.4
$ws.Range("K7").Value = 2313.4
$ws.Range("M7").Value = -2201.4
$ws.Range("H22").Value = 3975
$ws.Range("J22").Value = 2833.3333
$ws.Range("L22").Value = 2833.3333
$ws.Range("N22").Value = -3423.3333
$ws.Range("H27").Value = 3975
$ws.Range("J27").Value = 2833.3333
$ws.Range("L27").Value = 2833.3333
$ws.Range("N27").Value = -3047.3333
$ws.Range("H40").Value = 2985
$ws.Range("I40").Value = 2289.6
$ws.Range("K40").Value = 2289.6
$ws.Range("M40").Value = -2153.6
$ws.Range("H126").Value = 3593.5833
$ws.Range("I126").Value = 2313.4
$ws.Range("K126").Value = 6940.200000000001
$ws.Range("M126").Value = -4470.200000000001
$ws.Range("H133").Value = 76442
$ws.Range("J133").Value = 76442
$ws.Range("L133").Value = 76442
$ws.Range("N133").Value = -81502

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 3200
$ws.Range("I29").Value = 3200
$ws.Range("K29").Value = 3200
$ws.Range("M29").Value = -2910
$ws.Range("H119").Value = 29996
$ws.Range("J119").Value = 29996
$ws.Range("L119").Value = 29996
$ws.Range("N119").Value = -39672
$ws.Range("H132").Value = 1654
$ws.Range("I132").Value = 1072.8
$ws.Range("K132").Value = 3218.4
$ws.Range("M132").Value = -688.3999999999996

Write-Output "Refreshed cached market values on ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR."